# Graviton tx_chain rework (Rev2B):
#  - swap saw1/mixer/saw2/amp2 cascade order and fold the standalone saw2
#    row into the main gain-budget table
#  - flip sign of the switch insertion-loss entry (zener installed backwards)
#  - add a parallel "meas" (measured) column block (H:N) that tracks the
#    budget numbers against a bench-measured cascade, plus a diff column (Q)
#  - relabel the F column header as "budget" and retarget the subtotal/total
#    formulas to the new row layout

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tx_chain")

# ---- helper: apply the common "0.0" numeric display used throughout ----
function Set1Dec($rng) {
    $rng.NumberFormat = "0.0"
}

# ===================== Row 10 (headers) =====================
$ws.Range("F10").Value = "budget"

$ws.Range("H10").Value = "expected"
$ws.Range("K10").Value = "meas"
$ws.Range("M10").Value = "expected"
$ws.Range("Q10").Value = "diff"

# ===================== Row 11 =====================
# E11/F11/G11 unchanged (DAC out / mixer in (dBm), 1, dBm)
$ws.Range("H11").Formula = "=F11"
Set1Dec($ws.Range("H11"))

$ws.Range("K11").Formula = "=-32.9+30+3"
Set1Dec($ws.Range("K11"))
$ws.Range("L11").Value = "dBm"

# ===================== Row 12 =====================
# E12/F12/G12 unchanged (mixer filter, -1.5, dB)
$ws.Range("H12").Formula = "=H11+F12"
Set1Dec($ws.Range("H12"))

Set1Dec($ws.Range("K12"))
$ws.Range("L12").Value = "dB"
$ws.Range("M12").Formula = "=K11+F12"
Set1Dec($ws.Range("M12"))

# ===================== Row 13 =====================
# E13/F13/G13 unchanged (transformer, -0.35, dB)
$ws.Range("H13").Formula = "=H12+F13"
Set1Dec($ws.Range("H13"))

Set1Dec($ws.Range("K13"))
$ws.Range("L13").Value = "dB"
$ws.Range("M13").Formula = "=M12+F13"
Set1Dec($ws.Range("M13"))

# ===================== Row 14 =====================
# E14/F14/G14 unchanged (saw1 filter loss label "filter...", -7, dB)
$ws.Range("H14").Formula = "=H13+F14"
Set1Dec($ws.Range("H14"))

$ws.Range("K14").Formula = "=-38.3+30"
Set1Dec($ws.Range("K14"))
$ws.Range("L14").Value = "dB"
$ws.Range("M14").Formula = "=M13+F14"
Set1Dec($ws.Range("M14"))
Set1Dec($ws.Range("N14"))

# ===================== Row 15 =====================
# E15/F15/G15 unchanged (mixer filter label, -0.4, dB)
$ws.Range("H15").Formula = "=H14+F15"
Set1Dec($ws.Range("H15"))

Set1Dec($ws.Range("K15"))
$ws.Range("L15").Value = "dB"
$ws.Range("N15").Formula = "=K14+F15"
Set1Dec($ws.Range("N15"))

# ===================== Row 16 (was: mixer 24) -> now: saw1 -2 ===========
$ws.Range("E16").Value = "saw1 (dB), SF2098H"
$ws.Range("F16").Value = -2
Set1Dec($ws.Range("F16"))
# G16 unchanged ("dB")
$ws.Range("H16").Formula = "=H15+F16"
Set1Dec($ws.Range("H16"))

Set1Dec($ws.Range("K16"))
$ws.Range("L16").Value = "dB"
$ws.Range("N16").Formula = "=N15+F16"
Set1Dec($ws.Range("N16"))

# ===================== Row 17 (was: saw1 -2) -> now: mixer 24 ===========
$ws.Range("E17").Value = "mixer (dB), MAX2031"
$ws.Range("F17").Value = 24
Set1Dec($ws.Range("F17"))
# G17 unchanged ("dB")
$ws.Range("H17").Formula = "=H16+F17"
Set1Dec($ws.Range("H17"))

Set1Dec($ws.Range("K17"))
$ws.Range("L17").Value = "dB"
$ws.Range("N17").Formula = "=N16+F17"
Set1Dec($ws.Range("N17"))

# ===================== Row 18 (was: amp2 15.6) -> now: saw2 -2 ==========
$ws.Range("E18").Value = "saw2 (dB), SF2098H"
$ws.Range("F18").Value = -2
Set1Dec($ws.Range("F18"))
# G18 unchanged ("dB")
$ws.Range("H18").Formula = "=H17+F18"
Set1Dec($ws.Range("H18"))

Set1Dec($ws.Range("K18"))
$ws.Range("L18").Value = "dB"
$ws.Range("N18").Formula = "=N17+F18"
Set1Dec($ws.Range("N18"))

# ===================== Row 19 (new: amp2 15.6, folded in from old row22) ====
$ws.Range("E19").Value = "amp2 (dB), ALM-31122"
$ws.Range("F19").Value = 15.6
Set1Dec($ws.Range("F19"))
$ws.Range("G19").Value = "dB"
$ws.Range("H19").Formula = "=H18+F19"
Set1Dec($ws.Range("H19"))

Set1Dec($ws.Range("K19"))
$ws.Range("L19").Value = "dB"
$ws.Range("N19").Formula = "=N18+F19"
Set1Dec($ws.Range("N19"))

# ===================== Row 20 (new blank separator row) =====================
Set1Dec($ws.Range("F20"))
Set1Dec($ws.Range("K20"))

# ===================== Row 21 (subtotal "amp output", was row 20) ===========
$ws.Range("E21").Value = "amp output"
$ws.Range("F21").Formula = "=SUM(F11:F19)"
Set1Dec($ws.Range("F21"))
$ws.Range("G21").Value = "dBm"
Set1Dec($ws.Range("K21"))
$ws.Range("L21").Value = "dBm"

# ===================== Row 22 (blank separator, was standalone saw2 row) ====
Set1Dec($ws.Range("F22"))
Set1Dec($ws.Range("K22"))

# ===================== Row 23 (switch, sign flipped) =========================
# E23/G23 unchanged ("switch", "dB")
$ws.Range("F23").Value = -0.27
Set1Dec($ws.Range("F23"))
$ws.Range("H23").Formula = "=H19+F23"
Set1Dec($ws.Range("H23"))

$ws.Range("K23").Value = -0.27
Set1Dec($ws.Range("K23"))
$ws.Range("L23").Value = "dB"
$ws.Range("N23").Formula = "=N19+K23"
Set1Dec($ws.Range("N23"))

# ===================== Row 24 (blank separator) =====================
Set1Dec($ws.Range("F24"))
Set1Dec($ws.Range("K24"))

# ===================== Row 25 (subtotal "output power (dBm)") =====================
# E25/G25 unchanged
$ws.Range("F25").Formula = "=SUM(F21:F24)"
Set1Dec($ws.Range("F25"))
Set1Dec($ws.Range("K25"))
$ws.Range("L25").Value = "dBm"

# ===================== Row 26 (blank separator) =====================
Set1Dec($ws.Range("F26"))
Set1Dec($ws.Range("K26"))

# ===================== Row 27 (cable, estimate) =====================
# E27/G27 unchanged ("cable", "dB, estimated")
$ws.Range("H27").Formula = "=H23+F27"
Set1Dec($ws.Range("H27"))

$ws.Range("K27").Value = -0.5
Set1Dec($ws.Range("K27"))
$ws.Range("L27").Value = "dB, estimated"
$ws.Range("N27").Formula = "=N23+K27"
Set1Dec($ws.Range("N27"))

# ===================== Row 28 (thick bottom border separator) =====================
Set1Dec($ws.Range("K28"))

# ===================== Row 29 (expected / total) =====================
# E29/G29 unchanged ("expected", "dBm")
$ws.Range("F29").Formula = "=F25+F27"
Set1Dec($ws.Range("F29"))
Set1Dec($ws.Range("K29"))
$ws.Range("L29").Value = "dBm"

# ===================== Row 30 (actual / connector) =====================
# E30/F30/G30 unchanged ("actual", =-5.53+30, "dBm")
$ws.Range("K30").Value = 25
Set1Dec($ws.Range("K30"))
$ws.Range("L30").Value = "dBm"
$ws.Range("Q30").Formula = "=K30-N27"
Set1Dec($ws.Range("Q30"))

# ===================== Row 32 (difference) =====================
# E32/G32 unchanged ("difference", "dB")
# F32 formula unchanged ("=F29-F30") - value changes automatically w/ recalc
Set1Dec($ws.Range("K32"))
$ws.Range("L32").Value = "dB"

# ---- final touches: selection / active view (matches authored state) ----
$ws.Activate()
$ws.Range("H17").Select()

Write-Host "tx_chain rework applied"
